$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F header ("状态" = Status) ---
$ws.Range("F1").Value = "状态"

# --- Row-by-row updates: IsTrain (E) flips + new status (F) / note (G) columns ---
# Row 2
$ws.Range("F2").Value = "done"
# Row 3
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "缺少md"
# Row 4
$ws.Range("F4").Value = "done"
# Row 5
$ws.Range("F5").Value = "缺少md"
# Row 6
$ws.Range("F6").Value = "md飞了"
$ws.Range("G6").Value = "15026 frames"
# Row 7
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "缺少md"
# Row 8
$ws.Range("F8").Value = "md飞了"
$ws.Range("G8").Value = "10490 frames"
# Row 9
$ws.Range("F9").Value = "done"
# Row 10
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = "缺少md"
# Row 11
$ws.Range("F11").Value = "done"
# Row 12
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = "缺少md"
# Row 13
$ws.Range("F13").Value = "done"
# Row 14
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = "缺少md"
# Row 15
$ws.Range("F15").Value = "done"
# Row 16
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = "缺少md"
# Row 17
$ws.Range("F17").Value = "done"
# Row 18
$ws.Range("F18").Value = "done"
# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "done"
# Row 20
$ws.Range("F20").Value = "done"
# Row 21
$ws.Range("F21").Value = "done"
# Row 22
$ws.Range("F22").Value = "done"
# Row 23
$ws.Range("F23").Value = "done"
# Row 24
$ws.Range("F24").Value = "done"
# Row 25
$ws.Range("F25").Value = "done"

# --- AutoFilter over the full A1:F25 table ---
$ws.Range("A1:F25").AutoFilter()

# Register the hidden _FilterDatabase defined name Excel normally writes
# alongside an AutoFilter (sheet-scoped, hidden).
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$F`$25")
$fdb.Visible = $false

# --- Selection moved to E11 ---
$ws.Range("E11").Select()

# --- Page setup (A4 portrait) ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
